$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 2586
$ws1.Range("F5").Value = 913
$ws1.Range("F7").Value = 1878
$ws1.Range("F8").Value = 1731
$ws1.Range("F9").Value = 193
$ws1.Range("F11").Value = 2401
$ws1.Range("F12").Value = 508
$ws1.Range("F13").Value = 185
$ws1.Range("F18").Value = 8774
$ws1.Range("F20").Value = 6838
$ws1.Range("F21").Value = 11091
$ws1.Range("F26").Value = 531
$ws1.Range("F27").Value = 2441
$ws1.Range("F28").Value = 206
$ws1.Range("F30").Value = 2284
$ws1.Range("F31").Value = 418
$ws1.Range("F32").Value = 31
$ws1.Range("F33").Value = 4462
$ws1.Range("F34").Value = 627
$ws1.Range("F35").Value = 270
$ws1.Range("F36").Value = 20
$ws1.Range("F37").Value = 461

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F16").Value = 94

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 616
$ws3.Range("F5").Value = 82

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 616
$ws4.Range("F5").Value = 82
$ws4.Range("F7").Value = 2586
$ws4.Range("F9").Value = 913
$ws4.Range("F11").Value = 1878
$ws4.Range("F13").Value = 1731
$ws4.Range("F15").Value = 193
$ws4.Range("F17").Value = 508
$ws4.Range("F18").Value = 185
$ws4.Range("F23").Value = 8774
$ws4.Range("F25").Value = 6838
$ws4.Range("F26").Value = 11091
$ws4.Range("F33").Value = 531
$ws4.Range("F36").Value = 206
$ws4.Range("F37").Value = 31
$ws4.Range("F38").Value = 4462
$ws4.Range("F41").Value = 94
$ws4.Range("F45").Value = 461
